$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 451, pushing existing rows 451-478 down to 453-480.
$ws.Rows("451:452").Insert()

# Fill in the two new rows (451 and 452) with the new data records.
$ws.Range("A451").Value = 11
$ws.Range("B451").Value = "Vega Monumental Concepción"
$ws.Range("C451").Value = "Bíobío"
$ws.Range("D451").Value = 45013
$ws.Range("E451").Value = 8
$ws.Range("F451").Value = 100112002
$ws.Range("G451").Value = "Pimiento"
$ws.Range("H451").Value = "Morrón rojo"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 220
$ws.Range("K451").Value = 10000
$ws.Range("L451").Value = 11000
$ws.Range("M451").Value = 10545
$ws.Range("N451").Value = '$/caja 18 kilos'
$ws.Range("O451").Value = "Provincia de Limarí"
$ws.Range("P451").Value = 586
$ws.Range("Q451").Value = 18
$ws.Range("R451").Value = "Hortaliza"

$ws.Range("A452").Value = 11
$ws.Range("B452").Value = "Vega Monumental Concepción"
$ws.Range("C452").Value = "Bíobío"
$ws.Range("D452").Value = 45013
$ws.Range("E452").Value = 8
$ws.Range("F452").Value = 100112002
$ws.Range("G452").Value = "Pimiento"
$ws.Range("H452").Value = "Zafiro rojo"
$ws.Range("I452").Value = "Primera"
$ws.Range("J452").Value = 220
$ws.Range("K452").Value = 15000
$ws.Range("L452").Value = 16000
$ws.Range("M452").Value = 15455
$ws.Range("N452").Value = '$/caja 15 kilos'
$ws.Range("O452").Value = "Región de Arica y Parinacota"
$ws.Range("P452").Value = 1030
$ws.Range("Q452").Value = 15
$ws.Range("R452").Value = "Hortaliza"

# Ensure the date cells keep the same date number format as the rest of column D.
$ws.Range("D451").NumberFormat = $ws.Range("D453").NumberFormat()
$ws.Range("D452").NumberFormat = $ws.Range("D453").NumberFormat()
